$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update FechaSiniestro values in rows 2 and 3 (column H).
# Leading apostrophe keeps the cell as literal text (preserves quotePrefix style)
# instead of letting Excel auto-convert the string into a date value.
$ws.Range("H2").Value = "'19/03/2021"
$ws.Range("H3").Value = "'19/03/2021"

# Update NroPoliza values in rows 2 and 3 (column F)
$ws.Range("F2").Value = "04104015648"
$ws.Range("F3").Value = "04104015648"

# Update the sheet view - selection moves to C2, topLeftCell resets
$ws.Range("C2").Select()
